# Update row 2 values ("ADD results from server") across all sheets.
$wb = $excel.ActiveWorkbook

$sheet1Values = @(0, 487.6512737957041, 0, 0, 21227.96285102097, 0, 5667.147998863284, 0, 10965.327140292, 0, 0, 42857.36569267786, 8034.1746988965, 4555.555118045267, 5098.967661274243)

$otherSheetValues = @(889.3010627692065, 3966.462619977226, 0, 0, 42720.18362305129, 0, 5667.147998863284, 0, 33883.67254157657, 0, 0, 52955.41585787696, 17240.71739651479, 9549.220658140304, 9129.462057971876)

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    if ($sheetName -eq "2025") {
        $values = $sheet1Values
    } else {
        $values = $otherSheetValues
    }
    for ($col = 1; $col -le 15; $col++) {
        $ws.Cells.Item(2, $col).Value = $values[$col - 1]
    }
}
